# KAM Dashboard FY26 - change NPS achievement to improvement-based calculation.
# Adds a "Baseline (Q1)" column to the Annual KPIs sheet and updates the
# NPS Score row (Target, Achievement Till Date, Baseline).

$wb = $excel.ActiveWorkbook

# --- Annual KPIs sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("Annual KPIs")

# New header for the added baseline column.
$ws.Range("D3").Value = "Baseline (Q1)"

# NPS Score row: Target FY26 30 -> 19, Achievement Till Date -11 -> -12,
# and the new Baseline (Q1) value of -33.
$ws.Range("B6").Value = 19
$ws.Range("C6").Value = -12
$ws.Range("D6").Value = -33

# The header row used to be merged across A1:C1; now that the sheet has a
# fourth column it is left un-merged.
$ws.Range("A1:C1").UnMerge()
